$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data points (columns C and F) added to the "sample" sections
# (order matters for shared-string index assignment: "Gode målinger" must
# become string 22 and "Feilmålinger av transmittans" string 23)
$ws.Range("F6").Value = "Gode målinger"
$ws.Range("C6").Value = "Feilmålinger av transmittans"

$ws.Range("F12").Value = 80
$ws.Range("F13").Value = 98

$ws.Range("F30").Value = 82
$ws.Range("F31").Value = 100

$ws.Range("F39").Value = 92
$ws.Range("F40").Value = 99

# Update the view: scroll back to top and select D9 instead of F35
$ws.Range("D9").Select() | Out-Null
